# Updated cryptos list on Tue Feb 27 15:20:54 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for every tracked
# coin, and reflects the coinranking.com reshuffle that moved Stacks up
# (and FirstDigitalUSD / LidoDAOToken down) at rows 39-41, swapped
# ARBITRUM / NEARProtocol at rows 43-44, and replaced ThetaToken with
# ApeXProtocol at row 51.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text would otherwise be auto-detected as a
# number by Excel's Range.Value parser (e.g. "394.25", "1.00") --
# force them to Text format first so they round-trip as strings,
# matching the source data (which stores every Price/Volume cell as text).
$textCells = @(
    "D5",
    "D6",
    "D8",
    "D11",
    "D12",
    "D19",
    "D21",
    "D23",
    "D24",
    "D25",
    "D27",
    "D28",
    "D29",
    "D31",
    "D32",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D46",
    "D47",
    "D48",
    "D50",
    "D51"
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Row-by-row refresh of Price (D) / Volume 1h (E) columns, plus the
# coin-identity columns (B, C) for the rows that reshuffled position
# in the coinranking.com listing.
$ws.Range("D2").Value = '57.057.27'
$ws.Range("E2").Value = '  +9.63%  '
$ws.Range("D3").Value = '3.255.61'
$ws.Range("E3").Value = '  +4.37%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '394.25'
$ws.Range("E5").Value = '  -0.96%  '
$ws.Range("D6").Value = '109.25'
$ws.Range("E6").Value = '  +5.62%  '
$ws.Range("D7").Value = '3.250.34'
$ws.Range("E7").Value = '  +4.21%  '
$ws.Range("D8").Value = '0.568'
$ws.Range("E8").Value = '  +5.47%  '
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("E10").Value = '  +4.65%  '
$ws.Range("D11").Value = '39.29'
$ws.Range("E11").Value = '  +3.59%  '
$ws.Range("D12").Value = '0.0964'
$ws.Range("E12").Value = '  +12.11%  '
$ws.Range("E13").Value = '  +2.18%  '
$ws.Range("D14").Value = '3.766.00'
$ws.Range("E14").Value = '  +4.45%  '
$ws.Range("E15").Value = '  +5.19%  '
$ws.Range("E16").Value = '  +1.71%  '
$ws.Range("D17").Value = '3.252.26'
$ws.Range("E17").Value = '  +4.69%  '
$ws.Range("E18").Value = '  -0.91%  '
$ws.Range("D19").Value = '10.67'
$ws.Range("E19").Value = '  -1.04%  '
$ws.Range("D20").Value = '56.817.25'
$ws.Range("E20").Value = '  +9.41%  '
$ws.Range("D21").Value = '3.31'
$ws.Range("E21").Value = '  +2.51%  '
$ws.Range("E22").Value = '  +8.06%  '
$ws.Range("D23").Value = '13.05'
$ws.Range("E23").Value = '  +3.81%  '
$ws.Range("D24").Value = '300.07'
$ws.Range("E24").Value = '  +12.46%  '
$ws.Range("D25").Value = '74.32'
$ws.Range("E25").Value = '  +4.80%  '
$ws.Range("E26").Value = '  -3.12%  '
$ws.Range("D27").Value = '28.01'
$ws.Range("E27").Value = '  +2.48%  '
$ws.Range("D28").Value = '4.39'
$ws.Range("E28").Value = '  +4.20%  '
$ws.Range("D29").Value = '7.88'
$ws.Range("E29").Value = '  -1.62%  '
$ws.Range("E30").Value = '  +1.37%  '
$ws.Range("D31").Value = '7.22'
$ws.Range("E31").Value = '  -0.87%  '
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").Value = '  -0.05%  '
$ws.Range("E33").Value = '  +2.56%  '
$ws.Range("D34").Value = '11.00'
$ws.Range("E34").Value = '  +1.11%  '
$ws.Range("D35").Value = '38.18'
$ws.Range("E35").Value = '  +5.25%  '
$ws.Range("D36").Value = '0.0483'
$ws.Range("E36").Value = '  -1.91%  '
$ws.Range("D37").Value = '2.12'
$ws.Range("E37").Value = '  +1.63%  '
$ws.Range("D38").Value = '51.52'
$ws.Range("E38").Value = '  +3.01%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = '3.09'
$ws.Range("E39").Value = '  +17.60%  '
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  -0.22%  '
$ws.Range("B41").Value = 'LidoDAOToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D41").Value = '3.48'
$ws.Range("E41").Value = '  +1.98%  '
$ws.Range("D42").Value = '134.25'
$ws.Range("E42").Value = '  +2.86%  '
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").Value = '1.91'
$ws.Range("E43").Value = '  +1.92%  '
$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D44").Value = '4.01'
$ws.Range("E44").Value = '  -1.22%  '
$ws.Range("E45").Value = '  +3.58%  '
$ws.Range("D46").Value = '17.19'
$ws.Range("D47").Value = '0.284'
$ws.Range("E47").Value = '  -3.25%  '
$ws.Range("D48").Value = '21.93'
$ws.Range("E48").Value = '  -0.69%  '
$ws.Range("D49").Value = '2.147.88'
$ws.Range("E49").Value = '  +2.82%  '
$ws.Range("D50").Value = '2.07'
$ws.Range("E50").Value = '  +0.90%  '
$ws.Range("B51").Value = 'ApeXProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D51").Value = '2.38'
$ws.Range("E51").Value = '  -3.05%  '
